# Nexial io-showcase.xlsx edit:
# Add a new Web command "openIgnoreTimeout(url)" into the hidden '#system'
# lookup sheet's "web" column (column V), keeping the list alphabetically
# sorted. The new entry belongs right after "openHttpBasic(url,username,password)"
# (row 72) and before "refresh()" (row 73), so every existing entry from
# row 73 down to row 118 shifts down by one row (to 74..119), and the
# "web" named range grows from $V$2:$V$118 to $V$2:$V$119.
#
# Only column V is affected - column F (the unrelated "desktop" list that
# also lives on this sheet) must stay exactly where it is.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

$lastRow = 118
$insertRow = 73

# Shift existing column-V values down by one row, starting from the
# bottom so we don't clobber a value before it's been copied down.
for ($r = $lastRow; $r -ge $insertRow; $r--) {
    $srcAddr = "V" + $r
    $dstAddr = "V" + ($r + 1)
    $ws.Range($dstAddr).Value2 = $ws.Range($srcAddr).Value2
}

# Insert the new command in the now-vacated row.
$ws.Range("V73").Value2 = "openIgnoreTimeout(url)"

# Grow the "web" named range to cover the new row.
$names = $wb.Names
$webName = $names.Item("web")
$webName.RefersTo = "='#system'!`$V`$2:`$V`$119"
